# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.080.53'
$ws.Range('E2').Value = '  +0.46%  '
# Row 3
$ws.Range('D3').Value = '1.777.44'
$ws.Range('E3').Value = '  -0.29%  '
# Row 4
$ws.Range('E4').Value = '  +0.16%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '224.90'
$ws.Range('E5').Value = '  -0.61%  '
# Row 6
$ws.Range('E6').Value = '  -0.01%  '
# Row 7
$ws.Range('E7').Value = '  +0.23%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '31.56'
$ws.Range('E8').Value = '  -1.30%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.291'
$ws.Range('E9').Value = '  +0.21%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0683'
$ws.Range('E10').Value = '  +0.40%  '
# Row 11
$ws.Range('E11').Value = '  +1.03%  '
# Row 12
$ws.Range('D12').Value = '2.034.25'
$ws.Range('E12').Value = '  -0.29%  '
# Row 13
$ws.Range('D13').Value = '1.785.85'
$ws.Range('E13').Value = '  +0.01%  '
# Row 14
$ws.Range('E14').Value = '  -3.06%  '
# Row 15
$ws.Range('D15').Value = '34.080.37'
$ws.Range('E15').Value = '  +0.51%  '
# Row 16
$ws.Range('E16').Value = '  +0.74%  '
# Row 17
$ws.Range('E17').Value = '  +0.31%  '
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.41'
$ws.Range('E18').Value = '  -0.19%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.01'
$ws.Range('E19').Value = '  +0.89%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0784'
$ws.Range('E20').Value = '  +1.78%  '
# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.00'
$ws.Range('E21').Value = '  +3.59%  '
# Row 22
$ws.Range('E22').Value = '  +0.16%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.07'
$ws.Range('E23').Value = '  +0.17%  '
# Row 24
$ws.Range('E24').Value = '  -1.66%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '160.56'
$ws.Range('E25').Value = '  -0.68%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.10'
$ws.Range('E26').Value = '  -0.18%  '
# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.16'
$ws.Range('E27').Value = '  +0.04%  '
# Row 28
$ws.Range('E28').Value = '  +0.79%  '
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.30%  '
# Row 30
$ws.Range('E30').Value = '  -0.54%  '
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0515'
$ws.Range('E31').Value = '  +0.52%  '
# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.70'
$ws.Range('E32').Value = '  +2.38%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.68'
$ws.Range('E33').Value = '  +4.00%  '
# Row 34
$ws.Range('E34').Value = '  -2.54%  '
# Row 35
$ws.Range('D35').Value = '1.445.42'
$ws.Range('E35').Value = '  +3.65%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.653'
$ws.Range('E36').Value = '  +1.56%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.41'
$ws.Range('E37').Value = '  +4.65%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0189'
$ws.Range('E38').Value = '  +1.08%  '
# Row 39
$ws.Range('E39').Value = '  +0.26%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.37'
$ws.Range('E40').Value = '  +0.67%  '
# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '79.95'
$ws.Range('E41').Value = '  +0.41%  '
# Row 42
$ws.Range('E42').Value = '  +1.74%  '
# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.911'
$ws.Range('E43').Value = '  -0.49%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.53'
$ws.Range('E44').Value = '  -0.59%  '
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0515'
$ws.Range('E45').Value = '  +1.25%  '
# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.05'
$ws.Range('E46').Value = '  +2.70%  '
# Row 47
$ws.Range('E47').Value = '  -0.43%  '
# Row 48
$ws.Range('D48').Value = '1.935.71'
$ws.Range('E48').Value = '  -0.33%  '
# Row 49
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.29%  '
# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0130'
$ws.Range('E50').Value = '  -7.13%  '
# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '103.45'
$ws.Range('E51').Value = '  -3.09%  '

# Restore the default (General) style on cells that were temporarily
# forced to text format, so their cell style matches the original sheet.
$ws.Range('D5').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').Style = "Normal"
